$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "330.90"
Set-TextValue $ws.Range("E2") "6.71%"
Set-TextValue $ws.Range("G2") "13"
Set-TextValue $ws.Range("D3") "40.28"
Set-TextValue $ws.Range("E3") "8.30%"
Set-TextValue $ws.Range("G3") "13"
Set-TextValue $ws.Range("D4") "5.579"
Set-TextValue $ws.Range("E4") "8.80%"
Set-TextValue $ws.Range("G4") "13"
Set-TextValue $ws.Range("D5") "0.08118"
Set-TextValue $ws.Range("E5") "3.72%"
Set-TextValue $ws.Range("G5") "13"
Set-TextValue $ws.Range("D6") "8.682"
Set-TextValue $ws.Range("E6") "4.82%"
Set-TextValue $ws.Range("G6") "13"
Set-TextValue $ws.Range("D7") "1.970"
Set-TextValue $ws.Range("E7") "4.89%"
Set-TextValue $ws.Range("G7") "13"
Set-TextValue $ws.Range("E8") "-0.71%"
Set-TextValue $ws.Range("G8") "13"
Set-TextValue $ws.Range("D9") "0.9488"
Set-TextValue $ws.Range("E9") "2.85%"
Set-TextValue $ws.Range("G9") "13"
Set-TextValue $ws.Range("D10") "0.1260"
Set-TextValue $ws.Range("E10") "17.08%"
Set-TextValue $ws.Range("G10") "13"
Set-TextValue $ws.Range("E11") "4.42%"
Set-TextValue $ws.Range("G11") "13"
Set-TextValue $ws.Range("D12") "0.09198"
Set-TextValue $ws.Range("E12") "3.81%"
Set-TextValue $ws.Range("G12") "13"
Set-TextValue $ws.Range("D13") "0.03596"
Set-TextValue $ws.Range("E13") "8.39%"
Set-TextValue $ws.Range("G13") "13"
Set-TextValue $ws.Range("D14") "0.09596"
Set-TextValue $ws.Range("E14") "0.16%"
Set-TextValue $ws.Range("G14") "13"
Set-TextValue $ws.Range("D15") "0.001323"
Set-TextValue $ws.Range("E15") "-3.81%"
Set-TextValue $ws.Range("G15") "13"
Set-TextValue $ws.Range("D16") "0.006385"
Set-TextValue $ws.Range("E16") "10.20%"
Set-TextValue $ws.Range("G16") "13"
Set-TextValue $ws.Range("D17") "3.368"
Set-TextValue $ws.Range("E17") "-0.86%"
Set-TextValue $ws.Range("G17") "13"
Set-TextValue $ws.Range("D18") "4.550"
Set-TextValue $ws.Range("E18") "3.46%"
Set-TextValue $ws.Range("G18") "13"
Set-TextValue $ws.Range("D19") "0.3520"
Set-TextValue $ws.Range("E19") "2.56%"
Set-TextValue $ws.Range("G19") "13"
Set-TextValue $ws.Range("D20") "7.329"
Set-TextValue $ws.Range("E20") "15.64%"
Set-TextValue $ws.Range("G20") "13"
Set-TextValue $ws.Range("D21") "0.1342"
Set-TextValue $ws.Range("E21") "3.34%"
Set-TextValue $ws.Range("G21") "13"
Set-TextValue $ws.Range("E22") "1.76%"
Set-TextValue $ws.Range("G22") "13"
Set-TextValue $ws.Range("D23") "0.04430"
Set-TextValue $ws.Range("E23") "1.81%"
Set-TextValue $ws.Range("G23") "13"
Set-TextValue $ws.Range("D24") "0.001227"
Set-TextValue $ws.Range("E24") "2.16%"
Set-TextValue $ws.Range("G24") "13"
Set-TextValue $ws.Range("D25") "0.004347"
Set-TextValue $ws.Range("E25") "1.80%"
Set-TextValue $ws.Range("G25") "13"
Set-TextValue $ws.Range("D26") "0.0001201"
Set-TextValue $ws.Range("E26") "-14.27%"
Set-TextValue $ws.Range("G26") "13"
Set-TextValue $ws.Range("D27") "0.0003992"
Set-TextValue $ws.Range("E27") "37.58%"
Set-TextValue $ws.Range("G27") "13"
Set-TextValue $ws.Range("G28") "13"
Set-TextValue $ws.Range("G29") "13"
Set-TextValue $ws.Range("G30") "13"
Set-TextValue $ws.Range("G31") "13"
Set-TextValue $ws.Range("G32") "13"
Set-TextValue $ws.Range("G33") "13"
Set-TextValue $ws.Range("G34") "13"
Set-TextValue $ws.Range("G35") "13"
Set-TextValue $ws.Range("G36") "13"
Set-TextValue $ws.Range("G37") "13"
Set-TextValue $ws.Range("G38") "13"
Set-TextValue $ws.Range("D39") "0.02514"
Set-TextValue $ws.Range("E39") "16.65%"
Set-TextValue $ws.Range("G39") "13"
Set-TextValue $ws.Range("D40") "0.05243"
Set-TextValue $ws.Range("E40") "4.51%"
Set-TextValue $ws.Range("G40") "13"
Set-TextValue $ws.Range("D41") "0.007758"
Set-TextValue $ws.Range("E41") "2.35%"
Set-TextValue $ws.Range("G41") "13"
Set-TextValue $ws.Range("D42") "0.1433"
Set-TextValue $ws.Range("E42") "6.10%"
Set-TextValue $ws.Range("G42") "13"
Set-TextValue $ws.Range("D43") "0.008819"
Set-TextValue $ws.Range("E43") "3.64%"
Set-TextValue $ws.Range("G43") "13"
Set-TextValue $ws.Range("D44") "0.002142"
Set-TextValue $ws.Range("E44") "6.46%"
Set-TextValue $ws.Range("G44") "13"
Set-TextValue $ws.Range("D45") "0.01038"
Set-TextValue $ws.Range("E45") "28.06%"
Set-TextValue $ws.Range("G45") "13"
Set-TextValue $ws.Range("D46") "0.00006609"
Set-TextValue $ws.Range("E46") "0.93%"
Set-TextValue $ws.Range("G46") "13"
Set-TextValue $ws.Range("E47") "0.01%"
Set-TextValue $ws.Range("G47") "13"
Set-TextValue $ws.Range("D48") "0.002874"
Set-TextValue $ws.Range("E48") "-12.88%"
Set-TextValue $ws.Range("G48") "13"
Set-TextValue $ws.Range("D49") "0.002401"
Set-TextValue $ws.Range("E49") "66.32%"
Set-TextValue $ws.Range("G49") "13"
Set-TextValue $ws.Range("D50") "0.00002102"
Set-TextValue $ws.Range("E50") "0.01%"
Set-TextValue $ws.Range("G50") "13"
Set-TextValue $ws.Range("E51") "0.01%"
Set-TextValue $ws.Range("G51") "13"
